$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy number-format styling from row 19 down to the two new rows (20-21) ---
# Only the columns that actually carry a custom style in row 19 (G,H,O,P,Q,R,S,U,V);
# this avoids materialising stray empty cells in L/M/N/T like a full-row copy would.
$ws.Range("G19:H19").Copy()
$ws.Range("G20:H21").PasteSpecial(-4122)

$ws.Range("O19:S19").Copy()
$ws.Range("O20:S21").PasteSpecial(-4122)

$ws.Range("U19:V19").Copy()
$ws.Range("U20:V21").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Row 20 ---
$ws.Range("B20").Value = 17
$ws.Range("C20").Value = 286
$ws.Range("D20").Value = 159
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 45882.546527777777
$ws.Range("H20").Value = 45883.46597222222
$ws.Range("I20").Value = 7.45
$ws.Range("J20").Value = 5.2
$ws.Range("K20").Value = 2.5
$ws.Range("O20").Formula = "=H20-G20"
$ws.Range("P20").Formula = "=O20"
$ws.Range("Q20").Formula = "=I20-J20"
$ws.Range("R20").Formula = "=(F20-E20)/0.9982"
$ws.Range("S20").Formula = "=K20*P20"
$ws.Range("U20").Formula = "=Q20*1440/1324"
$ws.Range("V20").Formula = "=(1-ABS(U20-K20)/K20)*100%"

# --- Row 21 ---
$ws.Range("B21").Value = 18
$ws.Range("C21").Value = 286
$ws.Range("D21").Value = 159
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 45883.520833333336
$ws.Range("H21").Value = 45883.713194444441
$ws.Range("I21").Value = 5
$ws.Range("J21").Value = 4.55
$ws.Range("K21").Value = 2.5
$ws.Range("L21").Value = "Short Duration, less accuracy"
$ws.Range("O21").Formula = "=H21-G21"
$ws.Range("P21").Formula = "=O21"
$ws.Range("Q21").Formula = "=I21-J21"
$ws.Range("R21").Formula = "=(F21-E21)/0.9982"
$ws.Range("S21").Formula = "=K21*P21"
$ws.Range("U21").Formula = "=Q21*1440/277"
$ws.Range("V21").Formula = "=(1-ABS(U21-K21)/K21)*100%"

# --- Grow Table1 (Index..Note) to cover the two new rows ---
$lo1 = $ws.ListObjects.Item("Table1")
$lo1.Resize($ws.Range("B3:L21"))

# --- Conditional formatting (color scale) on the Accuracy column now spans through row 21 ---
$ws.Range("V4:V19").FormatConditions.Delete()
$ws.Range("V4:V21").FormatConditions.AddColorScale(3)

# --- Selection / view, matching the saved state in the source file ---
$ws.Range("L22").Select()

Write-Host "done"
